$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.072.30'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '''3.209.96'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''605.20'
$ws.Range('E5').Value = '  +4.98%  '
$ws.Range('D6').Value = '''154.22'
$ws.Range('E6').Value = '  +2.43%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '''3.208.74'
$ws.Range('E8').Value = '  +1.35%  '
$ws.Range('D9').Value = '''0.535'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('E10').Value = '  -0.39%  '
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '''0.508'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '''0.0000274'
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('D14').Value = '''38.73'
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').Value = '''3.734.73'
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = '''66.208.29'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').Value = '''7.46'
$ws.Range('E17').Value = '  +4.21%  '
$ws.Range('D18').Value = '''3.216.79'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = '''511.47'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').Value = '''15.60'
$ws.Range('E21').Value = '  +5.19%  '
$ws.Range('E22').Value = '  +0.66%  '
$ws.Range('D23').Value = '''15.25'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '''8.01'
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('D25').Value = '''85.30'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +3.39%  '
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('E30').Value = '  +9.90%  '
$ws.Range('E31').Value = '  +3.06%  '
$ws.Range('D32').Value = '''28.22'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').Value = '''6.64'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = '''0.0913'
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').Value = '''484.89'
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('E40').Value = '  -4.68%  '
$ws.Range('D41').Value = '''8.84'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').Value = '''2.51'
$ws.Range('E44').Value = '  +4.95%  '
$ws.Range('D45').Value = '''2.952.53'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('E46').Value = '  +6.33%  '
$ws.Range('D47').Value = '''28.87'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').Value = '''2.33'
$ws.Range('E50').Value = '  +3.23%  '
$ws.Range('D51').Value = '''34.26'
$ws.Range('E51').Value = '  +6.35%  '
